$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the "Swap promedio cámara 2021 - Diaria" sheet
# Columns: A=Serie(date), B..M = SPC values
$newRows = @(
    @{ Row = 184; A = "20-09-2021"; B = 3.85; C = 4.09; D = 4.26; E = 4.41; F = 4.86; G = -0.95; H = -0.02; I = 0.33; J = 0.72; K = 0.97; L = 1.49; M = 1.64 },
    @{ Row = 185; A = "21-09-2021"; B = 3.86; C = 4.08; D = 4.25; E = 4.4;  F = 4.84; G = -0.9;  H = 0.04;  I = 0.37; J = 0.75; K = 1.01; L = 1.5;  M = 1.66 },
    @{ Row = 186; A = "22-09-2021"; B = 3.92; C = 4.13; D = 4.3;  E = 4.46; F = 4.88; G = -0.86; H = 0.07000000000000001; I = 0.4;  J = 0.78; K = 1.04; L = 1.53; M = 1.7 },
    @{ Row = 187; A = "23-09-2021"; B = 4.04; C = 4.27; D = 4.44; E = 4.6;  F = 5.04; G = -0.92; H = 0.2;  I = 0.55; J = 0.89; K = 1.17; L = 1.68; M = 1.86 },
    @{ Row = 188; A = "24-09-2021"; B = 4.13; C = 4.37; D = 4.53; E = 4.69; F = 5.12; G = -1;    H = 0.25; I = 0.6;  J = 0.96; K = 1.23; L = 1.77; M = 1.94 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}
